# "Add files via upload" — rework Sheet1:
#   * insert a new "Area" column at G (pushing the old polygon_validation
#     column to H), and add a new blank "validation_remark" column at I
#   * re-type Farmer_ID / Field ID (columns A & B) as real numbers instead
#     of text
#   * populate the new Area column with the computed polygon areas, and
#     clear out the old "Good Polygon" / "Bad not fixable" remarks (G) —
#     only row 16 keeps a remark, moved over to H as "Modified"
#   * row 10's Coordinates (F) were re-digitised, so refresh that text too

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -------------------------------------------------------
# Clone G1's style onto the two new header cells before touching any
# values, so H1/I1 end up bold + bordered like the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = "Area"
$ws.Range("I1").Value = "validation_remark"

# Materialise blank (but present) H/I cells for every data row up front
# -- a plain "" assignment would just delete the cell again, so instead
# copy an unstyled, empty cell's formatting over the target range. That
# is enough to make Excel keep an empty <c> node without giving it any
# number format / style of its own.
$ws.Range("A2").Copy()
$ws.Range("H2:I19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- per-row data -------------------------------------------------------
$areas = @{
    2  = 19.23716898572448
    3  = 14.84156294075118
    4  = 4.86905480040647
    5  = 4.969527732006752
    6  = 1.291003838590159
    7  = 7.077876901175668
    8  = 10.15844958127223
    9  = 13.23720915046712
    10 = 11.75066274553506
    11 = 16.39976752549545
    12 = 5.822398904355146
    13 = 14.8626122193352
    14 = 15.109446777607
    15 = 13.41925492986151
    16 = 24.23821727873494
    17 = 11.66686954406143
    18 = 35.2552433506886
    19 = 16.81185589804555
}

$remarks = @{
    16 = "Modified"
}

for ($r = 2; $r -le 19; $r++) {
    # A (Farmer_ID) / B (Field ID): text -> number, same value.
    $ws.Range("A$r").Value = $ws.Range("A$r").Value2
    $ws.Range("B$r").Value = $ws.Range("B$r").Value2

    # G: drop the old text remark, write the numeric area instead.
    $ws.Range("G$r").Value = $areas[$r]

    # H: new home of the old polygon_validation remark -- the blank
    # cell was already materialised above, so only rows that actually
    # carry a remark need a value written.
    if ($remarks.ContainsKey($r)) {
        $ws.Range("H$r").Value = $remarks[$r]
    }

    # I: brand-new validation_remark column, blank for every row (cell
    # already materialised above -- nothing further to do).
}

# Row 10's Coordinates were re-captured/re-digitised.
$ws.Range("F10").Value = "74.04089781286125,30.04983873106723,0 74.04128419100407,30.04981460359325,0 74.04128374277002,30.04879338707973,0 74.04261714474397,30.0488075723396,0 74.04260118772308,30.05046927665487,0 74.04261920660612,30.05180426097061,0 74.04169773128763,30.05179305257542,0 74.04169593835532,30.05144363573627,0 74.04089790250822,30.05143338738495,0 74.04089781286125,30.04983873106723,0 74.04089781286125,30.04983873106723,0 "
